# Insert a new data row at row 222 (pushing the existing rows 222:243 down
# to 223:244) and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 222:243 down by inserting a new row at 222.
$ws.Rows.Item(222).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(222, 1).Value  = 8
$ws.Cells.Item(222, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(222, 3).Value  = "Coquimbo"
$ws.Cells.Item(222, 4).Value  = 45132
$ws.Cells.Item(222, 5).Value  = 4
$ws.Cells.Item(222, 6).Value  = 100112044
$ws.Cells.Item(222, 7).Value  = "Perejil"
$ws.Cells.Item(222, 8).Value  = "Sin especificar"
$ws.Cells.Item(222, 9).Value  = "Primera"
$ws.Cells.Item(222, 10).Value = 2000
$ws.Cells.Item(222, 11).Value = 2500
$ws.Cells.Item(222, 12).Value = 3000
$ws.Cells.Item(222, 13).Value = 2750
$ws.Cells.Item(222, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(222, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(222, 16).Value = 1833
$ws.Cells.Item(222, 17).Value = 1.5
$ws.Cells.Item(222, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by
# the rest of the "Fecha" column.
$ws.Cells.Item(222, 4).NumberFormat = $ws.Cells.Item(223, 4).NumberFormat
